$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 is "Chapter 23: Application Programming Interface (API)" (col A).
# The note "revise it from objects tomorrow" in B25 is cleared and the
# green "DONE" highlight is reset to a plain white fill now that the
# API chapter has only just been started (per commit message).
$ws.Range("B25").ClearContents()
$ws.Range("B25").Interior.ThemeColor = [Microsoft.Office.Interop.Excel.XlThemeColor]::xlThemeColorLight1

# Move the view/selection to where work resumed.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D24").Select()
